$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 713.8823
$ws.Range("I107").Value = 477.25
$ws.Range("J107").Value = 4500
$ws.Range("K107").Value = 477.25
$ws.Range("L107").Value = 4500
$ws.Range("M107").Value = 1442.75
$ws.Range("N107").Value = -8340
$ws.Range("H137").Value = 3217.838
$ws.Range("J137").Value = 6938.1113
$ws.Range("L137").Value = 20814.3339
$ws.Range("N137").Value = -25914.3339
$ws.Range("H138").Value = 928744.7
$ws.Range("I138").Value = 1221
$ws.Range("K138").Value = 3663
$ws.Range("M138").Value = 1477
$ws.Range("H140").Value = 94195
$ws.Range("J140").Value = 94195
$ws.Range("L140").Value = 94195
$ws.Range("N140").Value = -104555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 49896.43
$ws.Range("J31").Value = 110426
$ws.Range("L31").Value = 110426
$ws.Range("N31").Value = -111014
$ws.Range("H45").Value = 2517.7856
$ws.Range("I45").Value = 2645
$ws.Range("K45").Value = 2645
$ws.Range("M45").Value = -2268
$ws.Range("H61").Value = 22778502
$ws.Range("I61").Value = 45458492
$ws.Range("K61").Value = 45458492
$ws.Range("M61").Value = -45458280
$ws.Range("H74").Value = 7820234.5
$ws.Range("I74").Value = 14707797
$ws.Range("K74").Value = 14707797
$ws.Range("M74").Value = -14706923
$ws.Range("H77").Value = 7820234.5
$ws.Range("I77").Value = 14707797
$ws.Range("K77").Value = 73538985
$ws.Range("M77").Value = -73534617
$ws.Range("H97").Value = 1418.1305
$ws.Range("I97").Value = 979.4706
$ws.Range("K97").Value = 979.4706
$ws.Range("M97").Value = -483.4706
$ws.Range("H101").Value = 76897.625
$ws.Range("J101").Value = 76897.625
$ws.Range("L101").Value = 76897.625
$ws.Range("N101").Value = -83387.625
$ws.Range("H102").Value = 24383.334
$ws.Range("I102").Value = 24383.334
$ws.Range("K102").Value = 24383.334
$ws.Range("M102").Value = -22761.334
$ws.Range("H122").Value = 3721.6667
$ws.Range("I122").Value = 2123.75
$ws.Range("K122").Value = 6371.25
$ws.Range("M122").Value = -3921.25
$ws.Range("H132").Value = 14531
$ws.Range("I132").Value = 2994.5
$ws.Range("K132").Value = 8983.5
$ws.Range("M132").Value = -6453.5
$ws.Range("H136").Value = 22778502
$ws.Range("I136").Value = 45458492
$ws.Range("K136").Value = 136375476
$ws.Range("M136").Value = -136372926

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1709.6666
$ws.Range("J64").Value = 630
$ws.Range("L64").Value = 630
$ws.Range("N64").Value = -1080
$ws.Range("H67").Value = 1709.6666
$ws.Range("J67").Value = 630
$ws.Range("L67").Value = 630
$ws.Range("N67").Value = -2190
$ws.Range("H102").Value = 88180.125
$ws.Range("J102").Value = 96977
$ws.Range("L102").Value = 96977
$ws.Range("N102").Value = -103467
$ws.Range("H107").Value = 1655.5
$ws.Range("I107").Value = 2049.6155
$ws.Range("J107").Value = 923.5714
$ws.Range("K107").Value = 2049.6155
$ws.Range("L107").Value = 923.5714
$ws.Range("M107").Value = -129.6154999999999
$ws.Range("N107").Value = -4763.5714
$ws.Range("H134").Value = 28446.875
$ws.Range("I134").Value = 3148.9722
$ws.Range("K134").Value = 9446.9166
$ws.Range("M134").Value = -6911.9166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1062170.6
$ws.Range("J31").Value = 2334573.5
$ws.Range("L31").Value = 2334573.5
$ws.Range("N31").Value = -2335163.5
$ws.Range("H34").Value = 1062170.6
$ws.Range("J34").Value = 2334573.5
$ws.Range("L34").Value = 2334573.5
$ws.Range("N34").Value = -2334977.5
$ws.Range("H105").Value = 1782.8636
$ws.Range("I105").Value = 1525.75
$ws.Range("K105").Value = 1525.75
$ws.Range("M105").Value = 221.25
$ws.Range("H134").Value = 347703.44
$ws.Range("I134").Value = 401484.38
$ws.Range("K134").Value = 1204453.14
$ws.Range("M134").Value = -1201918.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3065
$ws.Range("I109").Value = 3030
$ws.Range("K109").Value = 9090
$ws.Range("M109").Value = -8050

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 101000
$ws.Range("J15").Value = 101000
$ws.Range("L15").Value = 101000
$ws.Range("N15").Value = -101576
$ws.Range("H70").Value = 5757.8696
$ws.Range("I70").Value = 4068.8
$ws.Range("J70").Value = 7057.154
$ws.Range("K70").Value = 4068.8
$ws.Range("L70").Value = 7057.154
$ws.Range("M70").Value = -3798.8
$ws.Range("N70").Value = -7597.154
$ws.Range("H73").Value = 5757.8696
$ws.Range("I73").Value = 4068.8
$ws.Range("J73").Value = 7057.154
$ws.Range("K73").Value = 4068.8
$ws.Range("L73").Value = 7057.154
$ws.Range("M73").Value = -3132.8
$ws.Range("N73").Value = -8929.154
$ws.Range("H80").Value = 6732.087
$ws.Range("J80").Value = 10360.25
$ws.Range("L80").Value = 10360.25
$ws.Range("N80").Value = -12356.25
$ws.Range("H81").Value = 101000
$ws.Range("J81").Value = 101000
$ws.Range("L81").Value = 101000
$ws.Range("N81").Value = -102996
$ws.Range("H83").Value = 6732.087
$ws.Range("J83").Value = 10360.25
$ws.Range("L83").Value = 51801.25
$ws.Range("N83").Value = -61785.25
$ws.Range("H84").Value = 101000
$ws.Range("J84").Value = 101000
$ws.Range("L84").Value = 303000
$ws.Range("N84").Value = -312984
$ws.Range("H122").Value = 15885
$ws.Range("I122").Value = 14416.333
$ws.Range("K122").Value = 43248.999
$ws.Range("M122").Value = -40798.999
$ws.Range("H132").Value = 142859950
$ws.Range("I132").Value = 142859950
$ws.Range("K132").Value = 428579850
$ws.Range("M132").Value = -428577320

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1081.3077
$ws.Range("I22").Value = 1072.25
$ws.Range("K22").Value = 1072.25
$ws.Range("M22").Value = -777.25
$ws.Range("H27").Value = 1081.3077
$ws.Range("I27").Value = 1072.25
$ws.Range("K27").Value = 1072.25
$ws.Range("M27").Value = -965.25
$ws.Range("H68").Value = 1800
$ws.Range("I68").Value = 1533.3334
$ws.Range("K68").Value = 1533.3334
$ws.Range("M68").Value = -784.3334
$ws.Range("H70").Value = 56775.332
$ws.Range("J70").Value = 56775.332
$ws.Range("L70").Value = 56775.332
$ws.Range("N70").Value = -57315.332
$ws.Range("H71").Value = 1800
$ws.Range("I71").Value = 1533.3334
$ws.Range("K71").Value = 7666.666999999999
$ws.Range("M71").Value = -3922.666999999999
$ws.Range("H73").Value = 56775.332
$ws.Range("J73").Value = 56775.332
$ws.Range("L73").Value = 56775.332
$ws.Range("N73").Value = -58647.332
$ws.Range("H80").Value = 68564
$ws.Range("J80").Value = 68564
$ws.Range("L80").Value = 68564
$ws.Range("N80").Value = -70810
$ws.Range("H82").Value = 1939.4667
$ws.Range("I82").Value = 1349.625
$ws.Range("J82").Value = 2613.5715
$ws.Range("K82").Value = 1349.625
$ws.Range("L82").Value = 2613.5715
$ws.Range("M82").Value = -988.625
$ws.Range("N82").Value = -3335.5715
$ws.Range("H83").Value = 68564
$ws.Range("J83").Value = 68564
$ws.Range("L83").Value = 205692
$ws.Range("N83").Value = -216924
$ws.Range("H85").Value = 1939.4667
$ws.Range("I85").Value = 1349.625
$ws.Range("J85").Value = 2613.5715
$ws.Range("K85").Value = 1349.625
$ws.Range("L85").Value = 2613.5715
$ws.Range("M85").Value = -101.625
$ws.Range("N85").Value = -5109.5715
$ws.Range("H100").Value = 3234.375
$ws.Range("I100").Value = 3411.5386
$ws.Range("K100").Value = 3411.5386
$ws.Range("M100").Value = -2870.5386
$ws.Range("H118").Value = 135404.5
$ws.Range("J118").Value = 135404.5
$ws.Range("L118").Value = 135404.5
$ws.Range("N118").Value = -138718.5
$ws.Range("H131").Value = 69999.664
$ws.Range("J131").Value = 69999.664
$ws.Range("L131").Value = 69999.664
$ws.Range("N131").Value = -80079.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 78997.664
$ws.Range("J140").Value = 78997.664
$ws.Range("L140").Value = 78997.664
$ws.Range("N140").Value = -89357.664
$ws.Range("H141").Value = 59993
$ws.Range("J141").Value = 59993
$ws.Range("L141").Value = 59993
$ws.Range("N141").Value = -70353
